$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder three country pairs (rank swapped since one overtook the other) ---
# Rows 46/47: China <-> Polonia
$ws.Range("A46").Value = "Polonia"
$ws.Range("A47").Value = "China"

# Rows 103/104: Tayikistan <-> Finlandia
$ws.Range("A103").Value = "Finlandia"
$ws.Range("A104").Value = "Tayikistan"

# Rows 205/206: Santa Lucia <-> Timor Oriental
$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("A206").Value = "Santa Lucia"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 12:00"

# --- Update numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 7: Rusia
$ws.Range("B7").Value = 1143571
$ws.Range("C7").Value = 7523
$ws.Range("D7").Value = 940150
$ws.Range("E7").Value = 183196
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 169
$ws.Range("H7").Value = 20225

# Row 18: Banglades
$ws.Range("B18").Value = 357873
$ws.Range("C18").Value = 1106
$ws.Range("D18").Value = 268777
$ws.Range("E18").Value = 83967
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = 5129

# Row 46: now Polonia (new data for the row that moved ahead)
$ws.Range("B46").Value = 85980
$ws.Range("C46").Value = 1584
$ws.Range("D46").Value = 67326
$ws.Range("E46").Value = 16230
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 32
$ws.Range("H46").Value = 2424

# Row 47: now China
$ws.Range("B47").Value = 85337
$ws.Range("C47").Value = 15
$ws.Range("D47").Value = 80536
$ws.Range("E47").Value = 167
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 4634

# Row 67: Austria
$ws.Range("B67").Value = 42214
$ws.Range("C67").Value = 714
$ws.Range("D67").Value = 33154
$ws.Range("E67").Value = 8273
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 787

# Row 97: Malasia
$ws.Range("B97").Value = 10769
$ws.Range("C97").Value = 82
$ws.Range("D97").Value = 9785
$ws.Range("E97").Value = 851
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 133

# Row 103: now Finlandia
$ws.Range("B103").Value = 9577
$ws.Range("C103").Value = 93
$ws.Range("D103").Value = 7850
$ws.Range("E103").Value = 1384
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 343

# Row 104: now Tayikistan
$ws.Range("B104").Value = 9562
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 8341
$ws.Range("E104").Value = 1147
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 74

# Row 142: Sri Lanka
$ws.Range("D142").Value = 3186
$ws.Range("E142").Value = 146
